# Append the latest NZ vaccination-by-date rows (15 Nov 2021 - 21 Nov 2021)
# to Sheet1, matching the style used by the existing data rows, and update
# the saved selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: row, date-serial, column B value, column C value
$newRows = @(
    @(272, 44515, 7764, 13678),
    @(273, 44516, 6664, 14367),
    @(274, 44517, 6197, 14211),
    @(275, 44518, 6635, 14028),
    @(276, 44519, 6833, 15775),
    @(277, 44520, 6002, 15499),
    @(278, 44521, 3079, 6772)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Range("A$row").Value = $r[1]
    $ws.Range("B$row").Value = $r[2]
    $ws.Range("C$row").Value = $r[3]
}

# Column A on the previous data rows uses the dd/mm/yyyy date style (cell
# style index 3); copy that formatting down onto the newly added A cells so
# the new rows look exactly like the rest of the table.
$ws.Range("A271").Copy() | Out-Null
$ws.Range("A272:A278").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the saved selection to the new last cell, as in the authored file.
$ws.Range("C271").Select() | Out-Null
